# Generate Report for Handoff
# A new source file (ffff6dfc087d-691f-4dff-a546-a7ad217a7911.md) has been picked
# up for localization, and the existing file's generated report identifiers
# (GUID-based names, xlf hashes, handoff timestamps) were regenerated.

$wb = $excel.ActiveWorkbook

$oldMd      = "2b92aafd-072d-474c-bfb9-9b19fa30b1e0.md"
$newMd      = "b3f2c905-d9c1-4df5-a750-b79115f9b6f2.md"
$newMd2     = "ffff6dfc087d-691f-4dff-a546-a7ad217a7911.md"
$cfgName    = ".localization-config"

$zhXlfOld   = "2b92aafd-072d-474c-bfb9-9b19fa30b1e0.938bf64c4ee846b80da5f81dc77114f694c4c031.zh-cn.xlf"
$zhXlfNew   = "b3f2c905-d9c1-4df5-a750-b79115f9b6f2.d92e6aac30bebe30c940b1bce4db051d723d70f1.zh-cn.xlf"
$deXlfOld   = "2b92aafd-072d-474c-bfb9-9b19fa30b1e0.938bf64c4ee846b80da5f81dc77114f694c4c031.de-de.xlf"
$deXlfNew   = "b3f2c905-d9c1-4df5-a750-b79115f9b6f2.d92e6aac30bebe30c940b1bce4db051d723d70f1.de-de.xlf"

$zhTimeNew  = "2016-03-07 10:17:00"
$deTimeNew  = "2016-03-07 10:17:15"
$epoch      = "0001-01-01 00:00:00"

$readyStatus   = "Ready for handoff"
$ignoredStatus = "Not to be localized"
$includeReason = "Include"
$ignoredReason = "Ignored"

$mdBase   = "https://github.com/OpenLocalizationTest/oltest/blob/6a8c047fae66d8781e19c537c9989cea9c618da5/e2e/"
$cfgUrl   = "https://github.com/OpenLocalizationTest/oltest/blob/6a8c047fae66d8781e19c537c9989cea9c618da5/.localization-config"
$zhXlfBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/885e5b1fa9382de9b5006161ab304d5cf5c8eba6/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/"
$deXlfBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8161ed787937b613d1f0aa3f218ababc48a490d3/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/"

# ---------------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# Make room for the new file's row (push the config row down).
$wsOverview.Rows.Item(3).Insert()

$wsOverview.Range("A2").Value = $newMd
$wsOverview.Range("B2").Value = $readyStatus
$wsOverview.Range("C2").Value = $readyStatus

$wsOverview.Range("A3").Value = $newMd2
$wsOverview.Range("B3").Value = $readyStatus
$wsOverview.Range("C3").Value = $readyStatus

$wsOverview.Range("A4").Value = $cfgName
$wsOverview.Range("B4").Value = $ignoredStatus
$wsOverview.Range("C4").Value = $ignoredStatus

# Hyperlinks don't track row-shifts on their own in this host, so rebuild them.
$wsOverview.Range("A1").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), ($mdBase + $newMd), "", "", $newMd)
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), ($mdBase + $newMd2), "", "", $newMd2)
$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), $cfgUrl, "", "", $cfgName)

# ---------------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Rows.Item(3).Insert()

$wsZh.Range("A2").Value = $newMd
$wsZh.Range("B2").Value = $readyStatus
$wsZh.Range("C2").Value = $zhXlfNew
$wsZh.Range("D2").Value = $zhTimeNew
$wsZh.Range("G2").Value = $epoch
$wsZh.Range("H2").Value = $includeReason

$wsZh.Range("A3").Value = $newMd2
$wsZh.Range("B3").Value = $readyStatus
$wsZh.Range("C3").Value = $zhXlfNew
$wsZh.Range("D3").Value = $zhTimeNew
$wsZh.Range("G3").Value = $epoch
$wsZh.Range("H3").Value = $includeReason

$wsZh.Range("A4").Value = $cfgName
$wsZh.Range("B4").Value = $ignoredStatus
$wsZh.Range("D4").Value = $epoch
$wsZh.Range("G4").Value = $epoch
$wsZh.Range("H4").Value = $ignoredReason

$wsZh.Range("A1").Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), ($mdBase + $newMd), "", "", $newMd)
$wsZh.Hyperlinks.Add($wsZh.Range("C2"), ($zhXlfBase + $zhXlfNew), "", "", $zhXlfNew)
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), ($mdBase + $newMd2), "", "", $newMd2)
$wsZh.Hyperlinks.Add($wsZh.Range("C3"), ($zhXlfBase + $zhXlfNew), "", "", $zhXlfNew)
$wsZh.Hyperlinks.Add($wsZh.Range("A4"), $cfgUrl, "", "", $cfgName)

# ---------------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Rows.Item(3).Insert()

$wsDe.Range("A2").Value = $newMd
$wsDe.Range("B2").Value = $readyStatus
$wsDe.Range("C2").Value = $deXlfNew
$wsDe.Range("D2").Value = $deTimeNew
$wsDe.Range("G2").Value = $epoch
$wsDe.Range("H2").Value = $includeReason

$wsDe.Range("A3").Value = $newMd2
$wsDe.Range("B3").Value = $readyStatus
$wsDe.Range("C3").Value = $deXlfNew
$wsDe.Range("D3").Value = $deTimeNew
$wsDe.Range("G3").Value = $epoch
$wsDe.Range("H3").Value = $includeReason

$wsDe.Range("A4").Value = $cfgName
$wsDe.Range("B4").Value = $ignoredStatus
$wsDe.Range("D4").Value = $epoch
$wsDe.Range("G4").Value = $epoch
$wsDe.Range("H4").Value = $ignoredReason

$wsDe.Range("A1").Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), ($mdBase + $newMd), "", "", $newMd)
$wsDe.Hyperlinks.Add($wsDe.Range("C2"), ($deXlfBase + $deXlfNew), "", "", $deXlfNew)
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), ($mdBase + $newMd2), "", "", $newMd2)
$wsDe.Hyperlinks.Add($wsDe.Range("C3"), ($deXlfBase + $deXlfNew), "", "", $deXlfNew)
$wsDe.Hyperlinks.Add($wsDe.Range("A4"), $cfgUrl, "", "", $cfgName)

Write-Host "Report regenerated for handoff."
